# Update the "last updated" timestamp in the title cell (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 15:02"

# Refresh per-country COVID figures (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes). A handful of countries swap
# rank (and therefore row) because the data table is kept sorted by "Casos totales"
# descending, so for those rows we also rewrite column A (country name).

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 4171358
$ws.Cells.Item(4, 3).Value = 1367
$ws.Cells.Item(4, 4).Value = 1980432
$ws.Cells.Item(4, 5).Value = 2043562
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 31
$ws.Cells.Item(4, 8).Value = 147364

# Row 6: India
$ws.Cells.Item(6, 1).Value = "India"
$ws.Cells.Item(6, 2).Value = 1307191
$ws.Cells.Item(6, 3).Value = 19061
$ws.Cells.Item(6, 4).Value = 828514
$ws.Cells.Item(6, 5).Value = 447870
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 162
$ws.Cells.Item(6, 8).Value = 30807

# Row 16: Arabia Saudita
$ws.Cells.Item(16, 1).Value = "Arabia Saudita"
$ws.Cells.Item(16, 2).Value = 262772
$ws.Cells.Item(16, 3).Value = 2378
$ws.Cells.Item(16, 4).Value = 215731
$ws.Cells.Item(16, 5).Value = 44369
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 37
$ws.Cells.Item(16, 8).Value = 2672

# Row 25: Catar
$ws.Cells.Item(25, 1).Value = "Catar"
$ws.Cells.Item(25, 2).Value = 108638
$ws.Cells.Item(25, 3).Value = 394
$ws.Cells.Item(25, 4).Value = 105420
$ws.Cells.Item(25, 5).Value = 3054
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 164

# Row 26: Irak
$ws.Cells.Item(26, 1).Value = "Irak"
$ws.Cells.Item(26, 2).Value = 104711
$ws.Cells.Item(26, 3).Value = 2485
$ws.Cells.Item(26, 4).Value = 71268
$ws.Cells.Item(26, 5).Value = 29231
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 90
$ws.Cells.Item(26, 8).Value = 4212

# Row 39: Kuwait
$ws.Cells.Item(39, 1).Value = "Kuwait"
$ws.Cells.Item(39, 2).Value = 62625
$ws.Cells.Item(39, 3).Value = 753
$ws.Cells.Item(39, 4).Value = 52915
$ws.Cells.Item(39, 5).Value = 9285
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 4
$ws.Cells.Item(39, 8).Value = 425

# Row 76: Dinamarca
$ws.Cells.Item(76, 1).Value = "Dinamarca"
$ws.Cells.Item(76, 2).Value = 13438
$ws.Cells.Item(76, 3).Value = 48
$ws.Cells.Item(76, 4).Value = 12340
$ws.Cells.Item(76, 5).Value = 485
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 613

# Row 82: Republica de Macedonia
$ws.Cells.Item(82, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(82, 2).Value = 9797
$ws.Cells.Item(82, 3).Value = 128
$ws.Cells.Item(82, 4).Value = 5254
$ws.Cells.Item(82, 5).Value = 4092
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 6
$ws.Cells.Item(82, 8).Value = 451

# Row 83: Bosnia y Herzegovina
$ws.Cells.Item(83, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(83, 2).Value = 9767
$ws.Cells.Item(83, 3).Value = 305
$ws.Cells.Item(83, 4).Value = 4555
$ws.Cells.Item(83, 5).Value = 4932
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 6
$ws.Cells.Item(83, 8).Value = 280

# Row 98: Croacia
$ws.Cells.Item(98, 1).Value = "Croacia"
$ws.Cells.Item(98, 2).Value = 4715
$ws.Cells.Item(98, 3).Value = 81
$ws.Cells.Item(98, 4).Value = 3555
$ws.Cells.Item(98, 5).Value = 1032
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 128

# Row 129: Islandia
$ws.Cells.Item(129, 1).Value = "Islandia"
$ws.Cells.Item(129, 2).Value = 1843
$ws.Cells.Item(129, 3).Value = 2
$ws.Cells.Item(129, 4).Value = 1823
$ws.Cells.Item(129, 5).Value = 10
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 10

# Row 141: Liberia
$ws.Cells.Item(141, 1).Value = "Liberia"
$ws.Cells.Item(141, 2).Value = 1135
$ws.Cells.Item(141, 3).Value = 18
$ws.Cells.Item(141, 4).Value = 621
$ws.Cells.Item(141, 5).Value = 443
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 71

# Row 142: Jordania
$ws.Cells.Item(142, 1).Value = "Jordania"
$ws.Cells.Item(142, 2).Value = 1131
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 1035
$ws.Cells.Item(142, 5).Value = 85
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 11

# Row 143: Niger
$ws.Cells.Item(143, 1).Value = "Niger"
$ws.Cells.Item(143, 2).Value = 1124
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 1022
$ws.Cells.Item(143, 5).Value = 33
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 69

# Row 190: San Martin (Parte Holandesa)
$ws.Cells.Item(190, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(190, 2).Value = 84
$ws.Cells.Item(190, 3).Value = 3
$ws.Cells.Item(190, 4).Value = 63
$ws.Cells.Item(190, 5).Value = 6
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 15

# Row 210: Groenlandia
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 13
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 13
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Islas Malvinas
$ws.Cells.Item(211, 1).Value = "Islas Malvinas"
$ws.Cells.Item(211, 2).Value = 13
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 13
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0
